$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.520102666666667
$ws.Range("H2").Value = 13.560308
$ws.Range("I2").Value = 0.9927775608668273
$ws.Range("J2").Value = 0.9927775608668273
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 10.37357685866534
$ws.Range("R2").Value = 93.36219172798802
$ws.Range("S2").Value = 0.01571326643064251
$ws.Range("T2").Value = 0.01571326643064251
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.520102666666667
$ws.Range("H3").Value = 13.560308
$ws.Range("I3").Value = 0.9927775608668273
$ws.Range("J3").Value = 0.9927775608668273
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 504.4059889622951
$ws.Range("R3").Value = 4539.653900660656
$ws.Range("S3").Value = 0.7640436661107467
$ws.Range("T3").Value = 0.7640436661107469
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.520102666666667
$ws.Range("H4").Value = 13.560308
$ws.Range("I4").Value = 0.9927775608668273
$ws.Range("J4").Value = 0.9927775608668273
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 140.6318584470636
$ws.Range("R4").Value = 1265.686726023572
$ws.Range("S4").Value = 0.213020628325438
$ws.Range("T4").Value = 0.213020628325438
$ws.Range("I5").Value = 0.007222439133172593
$ws.Range("J5").Value = 0.007222439133172593
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 0.07546758751233335
$ws.Range("R5").Value = 0.679208287611
$ws.Range("S5").Value = 0.0001143137343671924
$ws.Range("T5").Value = 0.0001143137343671924
$ws.Range("I6").Value = 0.007222439133172593
$ws.Range("J6").Value = 0.007222439133172593
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("S6").Value = 0.005558404108925201
$ws.Range("T6").Value = 0.005558404108925202
$ws.Range("I7").Value = 0.007222439133172593
$ws.Range("J7").Value = 0.007222439133172593
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("R7").Value = 9.207848465459
$ws.Range("S7").Value = 0.001549721289880199
$ws.Range("T7").Value = 0.001549721289880199
